$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Merge the two intro runs ("For this PA..." / "While our current...")
#    A no-op Find/Replace spanning both runs causes Word to recombine
#    them into a single run once formatting is confirmed identical.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "calculate various mathematical formulas.  While our current",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "calculate various mathematical formulas.  While our current", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Replace the legacy HYPERLINK field (<cmath> rendered via field
#    codes) with a real w:hyperlink run.
# ------------------------------------------------------------------
$cmathField = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $fld = $d.Fields.Item($i)
    if ($fld.Code.Text -match "HYPERLINK") {
        $cmathField = $fld
        break
    }
}
if ($cmathField -ne $null) {
    $cmathField.Delete()
}

# Turn "<> library" back into "<cmath> library" as plain text first...
$d.Content.Find.Execute(
    "Note that common mathematical functions can be found in the <> library.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Note that common mathematical functions can be found in the <cmath> library.", 2) | Out-Null

# ...then convert the word "cmath" into a proper hyperlink.
$linkRange = $d.Content
$linkRange.Find.Execute("cmath", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($linkRange, "http://www.cplusplus.com/reference/cmath/", "", "", "cmath") | Out-Null

# ------------------------------------------------------------------
# 3. Merge the "You must submit your assignment through Canvas no
#    later than midnight on " runs into one.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "You must submit your assignment through Canvas no later than midnight on",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "You must submit your assignment through Canvas no later than midnight on", 2) | Out-Null

# ------------------------------------------------------------------
# 4. "100 points possible" -> "50 points possible", with the _GoBack
#    bookmark relocated to mark the edit point (matches Word's normal
#    behaviour of moving _GoBack to the last text-entry location).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Your assignment will be judged by the following criteria (100 points possible):",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Your assignment will be judged by the following criteria (50 points possible):", 2) | Out-Null

$goBackRange = $d.Content
$goBackRange.Find.Execute("50", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
